$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set Runmode column (C) values for rows 3-7 to "Y" (same as row 2),
# which also drops the now-unused "N" shared string.
$ws.Range("C3:C7").Value = "Y"

# Move the active selection to C8.
$ws.Range("C8").Select()
